# Paragon workbook edit: handle case where we have no mapped outers
#
# 1. Delete the old "Sheet2" tab entirely.
# 2. Rename "Sheet3" to "unmapped-corp" and repoint its data/new column at
#    fake data describing an unmapped corp.
# 3. Make "unmapped-corp" the active tab and widen the window a bit.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Remove the old "Sheet2" worksheet entirely.
$wb.Worksheets.Item("Sheet2").Delete()

# 2. Repurpose "Sheet3" as the "unmapped-corp" example sheet.
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Name = "unmapped-corp"

# All CORP_CD values (column B, rows 2-25) become the fake unmapped code.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 2).Value = "F123D"
}

# Add a new column P explaining the test data, matching the header style.
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value = "this data is for testing a corp that we do not have mapped in outerskey"

# Selection/active-cell bookkeeping on the sheet itself.
$ws.Range("R16").Select()

# 3. Make "unmapped-corp" the active sheet/tab and resize the window.
$ws.Activate()
$excel.ActiveWindow.WindowState = -4143
$wb.Windows.Item(1).Width = 35700
